$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "uhdadhaksd"
$ws.Range("G1").Value = "dadasds"
$ws.Range("G4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G8").Value = 1
